# Add a PROFIT column (K) to the SalesData sheet, fix two mis-tagged
# product-type cells (INTERNALSTORAGE -> EXTERNALSTORAGE), and extend the
# AutoFilter / _FilterDatabase range to cover the new data extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SalesData")

# --- Fix the two rows that should have been EXTERNALSTORAGE all along ---
$ws.Range("E162").Value = "EXTERNALSTORAGE"
$ws.Range("E169").Value = "EXTERNALSTORAGE"

# --- New PROFIT header in K1 ---
$ws.Range("K1").Value = "PROFIT"

# --- PROFIT formula for K2:K183 (relative refs auto-adjust per row) ---
$profitFormula = '=IF(E2="INTERNALSTORAGE",H2*52%,IF(E2="SERVICE",H2*100%,IF(E2="DISPLAY",H2*50%,IF(E2="RAM",H2*65%,IF(E2="SOFTWARE",H2*100%,IF(E2="BATTERY",H2*70%,IF(E2="KEYBOARD",H2*60%,IF(E2="ADAPTOR",H2*45%,IF(E2="MAINBOARD",H2*15%,IF(E2="ACCESSORIES",H2*25%,IF(E2="VGA",H2*15%,IF(E2="POWERSUPPLY",H2*40%,IF(E2="PROCESSOR",H2*20%,IF(E2="SECOND",H2*150%,IF(E2="CASING",H2*15%,IF(E2="MONITOR",H2*5%,IF(E2="EXTERNALSTORAGE",100000,H2)))))))))))))))))'
$ws.Range("K2:K183").Formula = $profitFormula

# --- Extend the AutoFilter over the whole table, including new column K ---
$ws.AutoFilterMode = $false
$ws.Range("A1:K183").AutoFilter() | Out-Null

# --- Keep the _xlnm._FilterDatabase defined name in sync with the filter ---
$names = $wb.Names
for ($i = 1; $i -le $names.Count(); $i++) {
    $n = $names.Item($i)
    if ($n.Name() -eq "SalesData!_FilterDatabase") {
        $n.RefersTo = "=SalesData!`$A`$1:`$K`$183"
    }
}

# --- Match the saved view state: scrolled to show column B first, with
#     K2:K183 selected (active cell K2) ---
$win = $excel.ActiveWindow()
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("K2:K183").Select() | Out-Null
